$wb = $excel.ActiveWorkbook

# --- Sheet "List" (first sheet) ---
# Drop the "role.id" column (A keeps the role_roleName header/value,
# the old B column with role_roleName/roleName is removed).
$wsList = $wb.Worksheets.Item(1)
$wsList.Range("A1").Value = '${msg.getProperty(''role_roleName'')}'
$wsList.Range("A2").Value = '${role.roleName}'
$wsList.Range("B1:B2").ClearContents()

# --- Sheet "Search" (second sheet) ---
# Row 5 used to hold the role_roleName/roleName pair; it now becomes the
# role_id/id pair, and a brand-new row 6 carries the role_roleName/roleName
# pair that used to live in row 5.
$wsSearch = $wb.Worksheets.Item(2)
$wsSearch.Range("A5").Value = '${msg.getProperty(''role_id'')}'
$wsSearch.Range("B5").Value = '${id}'
$wsSearch.Range("A6").Value = '${msg.getProperty(''role_roleName'')}'
$wsSearch.Range("B6").Value = '${roleName}'
